$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new column E, styled like the other header cells (B1:D1):
# copy D1's format onto E1, then set its value.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "district"

# Data rows: first place name extracted from each row's District list in column D
$districts = @(
    "ADMIRALTY",
    "GREEN ISLAND",
    "CYBER PORT",
    "CHUNG HOM KOK",
    "SHEK O",
    "LEI YUE MUN",
    "SHAM TSENG",
    "SIU LAM",
    "LUNG KWU TAN",
    "LAM TEI",
    "MAI PO",
    "CLOSED AREA",
    "FEI NGO SHAN",
    "KAU SAI CHAU",
    "CLEAR WATER BAY",
    "CHI MA WAN",
    "HEI LING CHAU",
    "LAMMA ISLAND"
)

for ($i = 0; $i -lt $districts.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $districts[$i]
}
